# 03_PaschenuvZakon/data/PaschenuvZakon.xlsx
# "Small corrections, remade table 1"
#
# Adds a new "table 1" below the existing Paschen curve fit results on
# Sheet1: a small summary of the two measurement series (constant
# pressure / constant electrode distance) followed by a unit-conversion
# reference table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Table: B / C' / Uz / pd summary -------------------------------
$ws.Range("B47").Value = "B"
$ws.Range("C47").Value = "C'"
$ws.Range("D47").Value = "Uz"
$ws.Range("E47").Value = "pd"

$ws.Range("A48").Value = "Měření s kontantním tlakem"
$ws.Range("A48").WrapText = $False
$ws.Range("B48").Value = "310 30"
$ws.Range("C48").Value = "1.05 0.02"
$ws.Range("D48").Value = 291.6
$ws.Range("D48").NumberFormat = "0.0"
$ws.Range("E48").Value = 0.95

$ws.Range("A49").Value = "Měření s konstantní vzdáleností elektrod"
$ws.Range("A49").WrapText = $False
$ws.Range("B49").Value = "290 30"
$ws.Range("C49").Value = "0.79 0.01"
$ws.Range("D49").Value = 351.4
$ws.Range("D49").NumberFormat = "0.0"
$ws.Range("E49").Value = 1.23

# trailing styled-but-empty row, matching the formatting carried down
# from the rows above
$ws.Range("A50").WrapText = $False
$ws.Range("D50").NumberFormat = "0.0"

# ---- Unit-conversion reference table --------------------------------
$ws.Range("A51").Value = "Převody jednotek"

$ws.Range("A52").Value = "m2C-1"
$ws.Range("B52").Value = "m2A-1s-1"

$ws.Range("B53").Value = "Pa"
$ws.Range("C53").Value = "kgm-1s-2"

$ws.Range("B54").Value = "V"
$ws.Range("C54").Value = "m2kgs-3A-1"

$ws.Range("A55").Value = "VPa-1m-1"
$ws.Range("B55").Value = "m2s-1A-1"

# ---- Leave the view parked near the new table, like the author did --
[void]$ws.Range("C51").Select()
